# Update the "Metadata" worksheet (sheet 1) of the ValueSet workbook:
#  - bump Version to 0.1.7
#  - change Status from active -> draft
#  - update the Date
#  - update Contact info (publisher contact + add a named contact)
#  - insert a new "Jurisdiction" row
#  - shift Description/Purpose/Copyright/Immutable rows down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- capture the current (pre-edit) values of the rows that need to move ---
$descA = $ws.Range("A12").Value2   # "Description"
$descB = $ws.Range("B12").Value2   # "RxNorm codes for Defibrotide"
$purpA = $ws.Range("A13").Value2   # "Purpose"
$purpB = $ws.Range("B13").Value2   # ""
$copyA = $ws.Range("A14").Value2   # "Copyright"
$copyB = $ws.Range("B14").Value2   # ""
$immA  = $ws.Range("A15").Value2   # "Immutable"
$immB  = $ws.Range("B15").Value2   # "BooleanType[null]"

# --- make room for the new "Jurisdiction" row by shifting rows 12-15 down to 13-16 ---
# Row 16 is new territory, so first copy formatting from an existing formatted row.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

$ws.Range("A16").Value = $immA
$ws.Range("B16").Value = $immB

$ws.Range("A15").Value = $copyA
$ws.Range("B15").Value = $copyB

$ws.Range("A14").Value = $purpA
$ws.Range("B14").Value = $purpB

$ws.Range("A13").Value = $descA
$ws.Range("B13").Value = $descB

# --- row 12 becomes the new "Jurisdiction" row ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# --- update Version (row 3) ---
$ws.Range("B3").Value = "0.1.7"

# --- update Status (row 6) ---
$ws.Range("B6").Value = "draft"

# --- update Date (row 8) ---
$ws.Range("B8").Value = "2024-11-22T12:33:30-06:00"

# --- update Contact rows (10 and 11) ---
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$wb.Save()
